$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp shown in the title cell (A1)
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 15 de Abril de 2020 a las 13:22"

# Refresh country case figures. A handful of countries overtook their
# neighbours in "Casos totales", so the sorted rows below shuffle slightly
# (Bielorrusia, Uzbekistan, Bosnia y Herzegovina and Malta each move up a few
# spots, pushing the rows in between down by one).
# Row 20: Austria
$ws.Cells.Item(20, 1).Value = "Austria"
$ws.Cells.Item(20, 2).Value = 14290
$ws.Cells.Item(20, 3).Value = 64
$ws.Cells.Item(20, 4).Value = 8098
$ws.Cells.Item(20, 5).Value = 5799
$ws.Cells.Item(20, 6).Value = 232
$ws.Cells.Item(20, 7).Value = 9
$ws.Cells.Item(20, 8).Value = 393

# Row 45: Bielorrusia
$ws.Cells.Item(45, 1).Value = "Bielorrusia"
$ws.Cells.Item(45, 2).Value = 3728
$ws.Cells.Item(45, 3).Value = 447
$ws.Cells.Item(45, 4).Value = 203
$ws.Cells.Item(45, 5).Value = 3489
$ws.Cells.Item(45, 6).Value = 68
$ws.Cells.Item(45, 7).Value = 3
$ws.Cells.Item(45, 8).Value = 36

# Row 46: Panama
$ws.Cells.Item(46, 1).Value = "Panama"
$ws.Cells.Item(46, 2).Value = 3574
$ws.Cells.Item(46, 3).Value = 0
$ws.Cells.Item(46, 4).Value = 72
$ws.Cells.Item(46, 5).Value = 3407
$ws.Cells.Item(46, 6).Value = 106
$ws.Cells.Item(46, 7).Value = 0
$ws.Cells.Item(46, 8).Value = 95

# Row 47: Catar
$ws.Cells.Item(47, 1).Value = "Catar"
$ws.Cells.Item(47, 2).Value = 3428
$ws.Cells.Item(47, 3).Value = 0
$ws.Cells.Item(47, 4).Value = 373
$ws.Cells.Item(47, 5).Value = 3048
$ws.Cells.Item(47, 6).Value = 37
$ws.Cells.Item(47, 7).Value = 0
$ws.Cells.Item(47, 8).Value = 7

# Row 48: Luxemburgo
$ws.Cells.Item(48, 1).Value = "Luxemburgo"
$ws.Cells.Item(48, 2).Value = 3307
$ws.Cells.Item(48, 3).Value = 0
$ws.Cells.Item(48, 4).Value = 500
$ws.Cells.Item(48, 5).Value = 2740
$ws.Cells.Item(48, 6).Value = 30
$ws.Cells.Item(48, 7).Value = 0
$ws.Cells.Item(48, 8).Value = 67

# Row 49: Republica Dominicana
$ws.Cells.Item(49, 1).Value = "Republica Dominicana"
$ws.Cells.Item(49, 2).Value = 3286
$ws.Cells.Item(49, 3).Value = 0
$ws.Cells.Item(49, 4).Value = 162
$ws.Cells.Item(49, 5).Value = 2941
$ws.Cells.Item(49, 6).Value = 143
$ws.Cells.Item(49, 7).Value = 0
$ws.Cells.Item(49, 8).Value = 183

# Row 69: Uzbekistan
$ws.Cells.Item(69, 1).Value = "Uzbekistan"
$ws.Cells.Item(69, 2).Value = 1275
$ws.Cells.Item(69, 3).Value = 110
$ws.Cells.Item(69, 4).Value = 99
$ws.Cells.Item(69, 5).Value = 1172
$ws.Cells.Item(69, 6).Value = 8
$ws.Cells.Item(69, 7).Value = 0
$ws.Cells.Item(69, 8).Value = 4

# Row 70: Kazajistan
$ws.Cells.Item(70, 1).Value = "Kazajistan"
$ws.Cells.Item(70, 2).Value = 1275
$ws.Cells.Item(70, 3).Value = 43
$ws.Cells.Item(70, 4).Value = 220
$ws.Cells.Item(70, 5).Value = 1040
$ws.Cells.Item(70, 6).Value = 20
$ws.Cells.Item(70, 7).Value = 1
$ws.Cells.Item(70, 8).Value = 15

# Row 71: Eslovenia
$ws.Cells.Item(71, 1).Value = "Eslovenia"
$ws.Cells.Item(71, 2).Value = 1248
$ws.Cells.Item(71, 3).Value = 28
$ws.Cells.Item(71, 4).Value = 165
$ws.Cells.Item(71, 5).Value = 1022
$ws.Cells.Item(71, 6).Value = 34
$ws.Cells.Item(71, 7).Value = 5
$ws.Cells.Item(71, 8).Value = 61

# Row 72: Banglades
$ws.Cells.Item(72, 1).Value = "Banglades"
$ws.Cells.Item(72, 2).Value = 1231
$ws.Cells.Item(72, 3).Value = 219
$ws.Cells.Item(72, 4).Value = 49
$ws.Cells.Item(72, 5).Value = 1132
$ws.Cells.Item(72, 6).Value = 1
$ws.Cells.Item(72, 7).Value = 4
$ws.Cells.Item(72, 8).Value = 50

# Row 75: Bosnia y Herzegovina
$ws.Cells.Item(75, 1).Value = "Bosnia y Herzegovina"
$ws.Cells.Item(75, 2).Value = 1110
$ws.Cells.Item(75, 3).Value = 27
$ws.Cells.Item(75, 4).Value = 253
$ws.Cells.Item(75, 5).Value = 816
$ws.Cells.Item(75, 6).Value = 4
$ws.Cells.Item(75, 7).Value = 1
$ws.Cells.Item(75, 8).Value = 41

# Row 76: Lituania
$ws.Cells.Item(76, 1).Value = "Lituania"
$ws.Cells.Item(76, 2).Value = 1091
$ws.Cells.Item(76, 3).Value = 21
$ws.Cells.Item(76, 4).Value = 138
$ws.Cells.Item(76, 5).Value = 924
$ws.Cells.Item(76, 6).Value = 14
$ws.Cells.Item(76, 7).Value = 0
$ws.Cells.Item(76, 8).Value = 29

# Row 100: Malta
$ws.Cells.Item(100, 1).Value = "Malta"
$ws.Cells.Item(100, 2).Value = 399
$ws.Cells.Item(100, 3).Value = 6
$ws.Cells.Item(100, 4).Value = 44
$ws.Cells.Item(100, 5).Value = 352
$ws.Cells.Item(100, 6).Value = 4
$ws.Cells.Item(100, 7).Value = 0
$ws.Cells.Item(100, 8).Value = 3

# Row 101: Bolivia
$ws.Cells.Item(101, 1).Value = "Bolivia"
$ws.Cells.Item(101, 2).Value = 397
$ws.Cells.Item(101, 3).Value = 43
$ws.Cells.Item(101, 4).Value = 7
$ws.Cells.Item(101, 5).Value = 362
$ws.Cells.Item(101, 6).Value = 3
$ws.Cells.Item(101, 7).Value = 0
$ws.Cells.Item(101, 8).Value = 28

# Row 102: Jordania
$ws.Cells.Item(102, 1).Value = "Jordania"
$ws.Cells.Item(102, 2).Value = 397
$ws.Cells.Item(102, 3).Value = 0
$ws.Cells.Item(102, 4).Value = 235
$ws.Cells.Item(102, 5).Value = 155
$ws.Cells.Item(102, 6).Value = 5
$ws.Cells.Item(102, 7).Value = 0
$ws.Cells.Item(102, 8).Value = 7

# Row 103: Taiwan
$ws.Cells.Item(103, 1).Value = "Taiwan"
$ws.Cells.Item(103, 2).Value = 395
$ws.Cells.Item(103, 3).Value = 2
$ws.Cells.Item(103, 4).Value = 137
$ws.Cells.Item(103, 5).Value = 252
$ws.Cells.Item(103, 6).Value = 0
$ws.Cells.Item(103, 7).Value = 0
$ws.Cells.Item(103, 8).Value = 6

# Row 128: Brunei
$ws.Cells.Item(128, 1).Value = "Brunei"
$ws.Cells.Item(128, 2).Value = 136
$ws.Cells.Item(128, 3).Value = 0
$ws.Cells.Item(128, 4).Value = 108
$ws.Cells.Item(128, 5).Value = 27
$ws.Cells.Item(128, 6).Value = 2
$ws.Cells.Item(128, 7).Value = 0
$ws.Cells.Item(128, 8).Value = 1

